$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

function Get-ParagraphByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# 1) Signature block: "Maycon de Souza Silva" -> "Thamires de Oliveira " + "Sipionato"
#    (split into two runs, "Sipionato" wrapped in proofErr spellStart/spellEnd)
$p1 = Get-ParagraphByText("Maycon de Souza Silva")
$xml1 = "<w:p $wNs w14:paraId='543F356D' w14:textId='07224BF8' w:rsidR='00472520' w:rsidRPr='00321942' w:rsidRDefault='00472520' w:rsidP='00DF5ABB'>" + `
        "<w:pPr><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:b/><w:bCs/><w:lang w:val='pt-BR'/></w:rPr></w:pPr>" + `
        "<w:r w:rsidRPr='00321942'><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:b/><w:bCs/><w:lang w:val='pt-BR'/></w:rPr><w:t xml:space='preserve'>Thamires de Oliveira </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:b/><w:bCs/><w:lang w:val='pt-BR'/></w:rPr><w:t>Sipionato</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "</w:p>"
$p1.Range.InsertXML($xml1)

# 2) "O certificado só é válido enquanto trabalha para a Vestas" - merge the two runs
#    (drop the proofErr wrap around "Vestas")
$p2 = Get-ParagraphByText("O certificado")
$xml2 = "<w:p $wNs w14:paraId='0B036244' w14:textId='77777777' w:rsidR='00EC55BC' w:rsidRDefault='00EC55BC' w:rsidP='00884F44'>" + `
        "<w:pPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='18'/><w:szCs w:val='20'/><w:lang w:val='pt-BR'/></w:rPr></w:pPr>" + `
        "<w:r w:rsidRPr='00DF5ABB'><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='18'/><w:szCs w:val='20'/><w:lang w:val='pt-BR'/></w:rPr><w:t>O certificado só é válido enquanto trabalha para a Vestas</w:t></w:r>" + `
        "</w:p>"
$p2.Range.InsertXML($xml2)

# 3) "**Este curso é um treinamento interno de Vestas ... ao redor do mundo." - merge all
#    runs into one, dropping both proofErr wraps around "Vestas"
$p3 = Get-ParagraphByText("Este curso")
$xml3 = "<w:p $wNs w14:paraId='27F3BC61' w14:textId='77777777' w:rsidR='00EC55BC' w:rsidRPr='00DF5ABB' w:rsidRDefault='00EC55BC' w:rsidP='00EC55BC'>" + `
        "<w:pPr><w:tabs><w:tab w:val='left' w:pos='6565'/></w:tabs><w:spacing w:after='0'/><w:ind w:left='102'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='18'/><w:szCs w:val='20'/><w:lang w:val='pt-BR'/></w:rPr></w:pPr>" + `
        "<w:r w:rsidRPr='00DF5ABB'><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='18'/><w:szCs w:val='20'/><w:lang w:val='pt-BR'/></w:rPr><w:t>**Este curso é um treinamento interno de Vestas que complementar à NR-1 - DISPOSIÇÕES GERAIS e GERENCIAMENTO DE RISCOS OCUPACIONAIS, referente às atividades específicas das turbinas operadas por Vestas ao redor do mundo.</w:t></w:r>" + `
        "</w:p>"
$p3.Range.InsertXML($xml3)

"done"
